$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ordering of Attribute/Type rows (rows 2-21), per target diff
$data = @(
    @("case", "str"),
    @("operation_end_time", "datetime"),
    @("case:concept:name", "str"),
    @("event_id", "str"),
    @("lifecycle:state", "str"),
    @("time:timestamp", "datetime"),
    @("human_workstation_green_button_pressed", "float"),
    @("parameters", "dict"),
    @("identifier:id", "str"),
    @("process_model_id", "str"),
    @("requested_service_url", "str"),
    @("concept:name", "str"),
    @("current_task", "str"),
    @("response_status_code", "float"),
    @("org:resource", "str"),
    @("unsatisfied_condition_description", "str"),
    @("lifecycle:transition", "str"),
    @("complete_service_time", "str"),
    @("SubProcessID", "str"),
    @("planned_operation_time", "str")
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row = $row + 1
}
